$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (File) corrections ---
$ws.Cells.Item(86, 1).Value = 'Neurological/Intrathecal policy RIE.pdf'
$ws.Cells.Item(87, 1).Value = 'Diabetes_and_Glucose/Intravenous Insulin Therapy (not for DKA or HHS).pdf'
$ws.Cells.Item(93, 1).Value = 'Infection_and_sepsis/Trip Out of Unit infection guidance.pdf'

# --- Column B (Review date) updates ---
$bRows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 53, 54, 55, 56, 57, 58, 59, 60, 61, 62, 63, 64, 65, 66, 67, 68, 69, 70, 71, 72, 73, 74, 75, 76, 77, 78, 79, 80, 81, 83, 85, 86, 89, 90, 91, 92, 93, 94, 95, 96, 97, 98, 99, 100, 101, 102, 103, 105, 106, 107, 108, 109, 111, 112, 113, 114, 115, 117, 118, 119, 120, 121, 122, 123, 128, 129, 130, 131, 132, 133, 134, 135, 136, 137, 138, 139, 140, 141, 142, 143, 144, 145, 146, 147, 148, 149, 150, 151, 152, 153, 154, 155, 156, 157, 158, 159, 160, 161, 162, 163, 164, 165, 166, 167, 168, 169, 170, 171, 172, 173, 174)
$bVals = @(
    '2011-03-20',
    '2011-06-20',
    '2011-07-20',
    '2011-12-20',
    '2014-09-20',
    '2014-10-20',
    '2014-11-20',
    '2014-12-20',
    '2014-12-20',
    '2015-03-20',
    '2016-01-20',
    '2016-04-20',
    '2016-08-20',
    '2017-01-20',
    '2017-08-20',
    '2017-09-20',
    '2017-09-20',
    '2017-09-20',
    '2017-10-20',
    '2017-10-20',
    '2017-11-20',
    '2018-05-20',
    '2018-07-20',
    '2019-05-20',
    '2019-05-20',
    '2019-06-20',
    '2019-06-20',
    '2019-06-20',
    '2019-08-20',
    '2019-09-20',
    '2020-01-20',
    '2020-03-20',
    '2020-03-20',
    '2020-03-20',
    '2020-03-20',
    '2020-03-20',
    '2020-03-20',
    '2020-04-20',
    '2020-05-20',
    '2020-07-20',
    '2020-08-20',
    '2020-11-20',
    '2021-02-20',
    '2021-05-20',
    '2021-05-20',
    '2021-05-20',
    '2021-06-20',
    '2021-06-20',
    '2021-06-20',
    '2021-06-20',
    '2021-06-20',
    '2021-06-20',
    '2021-06-20',
    '2021-06-20',
    '2021-06-20',
    '2021-09-20',
    '2021-12-20',
    '2021-12-20',
    '2022-01-20',
    '2022-01-20',
    '2022-03-20',
    '2022-03-20',
    '2022-04-20',
    '2022-05-20',
    '2022-05-20',
    '2022-06-20',
    '2022-06-20',
    '2022-06-20',
    '2022-06-20',
    '2022-07-20',
    '2022-07-20',
    '2022-07-20',
    '2022-08-20',
    '2022-09-20',
    '2022-10-20',
    '2022-10-20',
    '2022-11-20',
    '2022-11-20',
    '2022-12-20',
    '2022-12-20',
    '2023-01-20',
    '2023-02-20',
    '2023-02-20',
    '2023-04-20',
    '2023-04-20',
    '2023-04-20',
    '2023-04-20',
    '2023-05-20',
    '2023-05-20',
    '2023-05-20',
    '2023-05-20',
    '2023-05-20',
    '2023-05-20',
    '2023-05-20',
    '2023-05-20',
    '2023-05-20',
    '2023-05-20',
    '2023-05-20',
    '2023-06-20',
    '2023-06-20',
    '2023-06-20',
    '2023-06-20',
    '2023-06-20',
    '2023-07-20',
    '2023-07-20',
    '2023-07-20',
    '2023-07-20',
    '2023-07-20',
    '2023-08-20',
    '2023-08-20',
    '2023-08-20',
    '2023-09-20',
    '2023-09-20',
    '2023-10-20',
    '2023-10-20',
    '2024-01-20',
    '2024-02-20',
    '2024-02-20',
    '2024-03-20',
    '2024-03-20',
    '2024-03-20',
    '2024-03-20',
    '2024-03-20',
    '2024-04-20',
    '2024-04-20',
    '2024-05-20',
    '2024-06-20',
    '2024-07-20',
    '2024-07-20',
    '2024-07-20',
    '2024-07-20',
    '2024-07-20',
    '2024-08-20',
    '2024-08-20',
    '2024-10-20',
    '2024-10-20',
    '2025-01-20',
    '2025-01-20',
    '2025-02-20',
    '2025-02-20',
    '2025-02-20',
    '2025-03-20',
    '2025-03-20',
    '2025-04-20',
    '2025-04-20',
    '2025-05-20',
    '2025-05-20',
    '2025-05-20',
    '2025-05-20',
    '2025-05-20',
    '2025-06-20',
    '2025-07-20',
    '2025-07-20',
    '2025-08-20',
    '2025-10-20',
    '2025-10-20',
    '2026-01-20',
    '2026-03-20',
    '2026-03-20',
    '2026-03-20',
    '2027-01-20',
    '2027-02-20'
)
for ($i = 0; $i -lt $bRows.Length; $i++) {
    $c = $ws.Cells.Item($bRows[$i], 2)
    $c.NumberFormat = "@"
    $c.Value = $bVals[$i]
}
